$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Y" values to C2 and C4 (previously blank), matching the bordered/
# wrap-text style already used by C3.
$ws.Range("C2").Value = "Y"
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C2").WrapText = $true

$ws.Range("C4").Value = "Y"
$ws.Range("C4").Borders.LineStyle = 1
$ws.Range("C4").WrapText = $true

# Update the active selection to C4
$ws.Range("C4").Select()
